# Append the latest gold-price row (row 5) to the sheet, mirroring the
# Jenkins job that refreshes the GoldData report each day.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "26-09-2025"
$ws.Range("B5").Value = "The price of gold in India today is ₹11,488 per gram for 24 karat gold, ₹10,530 per gram for 22 karat gold and ₹8,616 per gram for 18 karat gold (also called 999 gold)."
